# Update verification matrix: add a new row documenting RandomForestFunc.R
# and its checks, pushing the existing RandomForestGroupVarImp.R row down
# one slot and giving it a "Data loading..." verification note.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing "RandomForestGroupVarImp.R" entry (row 19) down to row 20.
$ws.Range("A20").Value = "RandomForestGroupVarImp.R"

# New row 19: RandomForestFunc.R + its verification note.
$ws.Range("A19").Value = "RandomForestFunc.R"
$ws.Range("B19").Value = "Each aspect of the RF forest estimation was checked during construction of the function. Construction of the CART trees is done with a package function from randomForestSRC. Estimation of OOB error and variable importance are implemented in a standard fashion and results have been checked to show expected behaviour."

# Row 20 now gets its own verification note for RandomForestGroupVarImp.R.
$ws.Range("B20").Value = "Data loading and data transformations were checked manually for correctness. Result look reasonable, especially compared to pairwise method in GenFigSourcesOfVarPairWise.R."

# Reflect the user's on-screen selection after the edit.
$ws.Range("B21").Select()
